$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 5
